$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Pin Mapping table and schematic" row (row 14) as completed
$ws.Range("B14").Value = "completed"

# Add new task row for "LCD"
$ws.Range("A15").Value = "LCD"

# Update selection to reflect the new active cell (B15)
$ws.Range("B15").Select()
